# Applies the cryptos list update described in the commit message
# (price/volume refresh from GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/percentage updates (safe from numeric auto-conversion) ---
$ws.Range("D2").Value = "29.249.09"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.866.65"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("E5").Value = "  +3.11%  "
$ws.Range("E6").Value = "  +1.47%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").Value = "1.913.83"
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").Value = "29.278.32"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("E18").Value = "  +0.28%  "
$ws.Range("E19").Value = "  +0.53%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "2.110.16"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  +5.75%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("E31").Value = "  +2.53%  "
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("E33").Value = "  +0.90%  "
$ws.Range("E34").Value = "  +1.57%  "
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  +0.65%  "
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "1.174.13"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("E41").Value = "  -1.04%  "
$ws.Range("E42").Value = "  +1.98%  "
$ws.Range("E43").Value = "  +2.23%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  -1.11%  "
$ws.Range("D47").Value = "2.010.34"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("E50").Value = "  +7.03%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E51").Value = "  +1.95%  "

# --- Price values that look like plain numbers: force Text format so Excel
#     keeps them as strings (matching the original inlineStr cell type),
#     then reset the style back to Normal so no stray number format remains.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7237"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07836"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3088"
$ws.Range("D9").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7218"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.237"
$ws.Range("D14").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.869"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "243.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007812"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.955"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1604"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.05"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.963"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.342"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.399"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05204"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.936"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.184"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.677"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01855"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.701"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9031"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.115"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "72.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "101.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5286"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.778"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.892"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.305"
$ws.Range("D51").Style = "Normal"
